$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 57 (pushes old row 57 -> 58, 58 -> 59, etc.)
$ws.Rows("57:57").Insert()

# Row 56: the original single entry (08:15-12:00) is split into two entries.
# Update the end time of row 56 to the new (earlier) end time.
$ws.Cells.Item(56, 5).Value2 = 0.4236111111111111

# New row 57 holds the second half of the split entry.
$ws.Cells.Item(57, 1).Value2 = 2014
$ws.Cells.Item(57, 2).Value2 = 3
$ws.Cells.Item(57, 3).Value2 = 10
$ws.Cells.Item(57, 4).Value2 = 0.44791666666666669
$ws.Cells.Item(57, 5).Value2 = 0.5

# Recompute the "time spent" formulas for rows 56 and 57 (minutes, then hours).
$ws.Cells.Item(56, 6).Formula = "=(E56-D56)*24*60"
$ws.Cells.Item(56, 7).Formula = "=F56/60"
$ws.Cells.Item(57, 6).Formula = "=(E57-D57)*24*60"
$ws.Cells.Item(57, 7).Formula = "=F57/60"

# The formerly-blank row (old 57) is now row 58; leave it untouched/blank.

# Update the summary rows (old 58/59/60 are now 59/60/61) so the SUM range
# picks up the newly inserted row and the chained formulas recompute.
$ws.Cells.Item(59, 6).Formula = "=SUM(F2:F58)"
$ws.Cells.Item(60, 6).Formula = "=F59/60"
$ws.Cells.Item(61, 6).Formula = "=F60/38.5"

# Restore the selection to the top-left data cell, matching the new view state.
$ws.Range("G2").Select()
